$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.617.69'
$ws.Range('E2').Value = '  +2.29%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.676.64'
$ws.Range('E3').Value = '  +2.72%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.05'
$ws.Range('E5').Value = '  +2.49%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.530'
$ws.Range('E6').Value = '  +2.39%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.48'
$ws.Range('E8').Value = '  +3.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.265'
$ws.Range('E9').Value = '  +2.62%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0647'
$ws.Range('E10').Value = '  +6.34%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0905'
$ws.Range('E11').Value = '  -0.06%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.918.26'
$ws.Range('E12').Value = '  +2.81%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.678.61'
$ws.Range('E13').Value = '  +2.86%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.610'
$ws.Range('E14').Value = '  +8.39%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '10.01'
$ws.Range('E15').Value = '  +8.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.02'
$ws.Range('E16').Value = '  +4.48%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.649.32'
$ws.Range('E17').Value = '  +2.34%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.47'
$ws.Range('E18').Value = '  +3.59%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.63'
$ws.Range('E19').Value = '  +0.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').Value = '  +3.30%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('E22').Value = '  +3.06%  '

$ws.Range('E23').Value = '  +2.27%  '

$ws.Range('E24').Value = '  -0.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.90'
$ws.Range('E25').Value = '  +0.55%  '

$ws.Range('E26').Value = '  +3.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.84'
$ws.Range('E27').Value = '  +2.30%  '

$ws.Range('E28').Value = '  +1.64%  '

$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('E30').Value = '  +1.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +3.96%  '

$ws.Range('E32').Value = '  +2.37%  '

$ws.Range('E33').Value = '  +3.66%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.496.69'
$ws.Range('E34').Value = '  +4.68%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.79'
$ws.Range('E35').Value = '  +7.91%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '84.38'
$ws.Range('E36').Value = '  +11.40%  '

$ws.Range('E37').Value = '  -0.40%  '

$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.600'
$ws.Range('E38').Value = '  +8.31%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0179'
$ws.Range('E39').Value = '  +5.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.68'
$ws.Range('E40').Value = '  -2.48%  '

$ws.Range('E41').Value = '  +0.41%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.838'
$ws.Range('E42').Value = '  +1.09%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.00'
$ws.Range('E43').Value = '  +0.09%  '

$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0501'
$ws.Range('E44').Value = '  +2.46%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  -0.02%  '

$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.02'
$ws.Range('E46').Value = '  +0.75%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.55'
$ws.Range('E47').Value = '  +3.26%  '

$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '50.94'
$ws.Range('E48').Value = '  -1.96%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.811.48'
$ws.Range('E49').Value = '  +2.15%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '94.56'
$ws.Range('E50').Value = '  +4.67%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0117'
$ws.Range('E51').Value = '  +2.05%  '
